$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 9, shifting existing rows 9-12 down to 10-13
$ws.Rows.Item(9).Insert()

# Fill in the new row 9 with data
$ws.Range("A9").Value = "####1er Torneo Federativo - C.A.E. - Sub 23, Prejuveniles y sub 23 (28 de Febrero y 1 de Marzo) - Juniors (Domingo 1 de Marzo)"
$ws.Range("B9").Value = "Juveniles"
$ws.Range("C9").Value = "caballeros"
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = "Cejas, Pedro Gonzalo"
$ws.Range("F9").Value = 89
$ws.Range("G9").Value = ""
$ws.Range("H9").Value = 89
